# Apply updated numeric values to the "展览" and "全部类型" worksheets.
# Both sheets contain identical data layouts, so the same set of cell
# updates is applied to each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("G2").Value = 58
    $ws.Range("F6").Value = 12478
    $ws.Range("F7").Value = 12478
    $ws.Range("F10").Value = 501
    $ws.Range("F12").Value = 1142
    $ws.Range("F14").Value = 13633
    $ws.Range("F15").Value = 13900
    $ws.Range("F20").Value = 1042
    $ws.Range("F23").Value = 406
    $ws.Range("F24").Value = 5001
    $ws.Range("F25").Value = 241
}
